# Add github action flow
# Adds a "Summary" worksheet (after "TestResults") that aggregates the daily
# pass/fail counts from the TestResults sheet.

$wb = $excel.ActiveWorkbook
$wsResults = $wb.Worksheets.Item(1)

# Insert the new sheet right after TestResults and name it "Summary"
$wsSummary = $wb.Worksheets.Add($null, $wsResults)
$wsSummary.Name = "Summary"

# Header row
$wsSummary.Range("A1").Value = "Date"
$wsSummary.Range("B1").Value = "Passed"
$wsSummary.Range("C1").Value = "Failed"
$wsSummary.Range("D1").Value = "Total"
$wsSummary.Range("E1").Value = "% Passed"

# Day 1 summary (16-05-2025)
$wsSummary.Range("A2").Value = "16-05-2025"
$wsSummary.Range("B2").Value = 10
$wsSummary.Range("C2").Value = 1
$wsSummary.Range("D2").Value = 11
$wsSummary.Range("E2").NumberFormat = "@"
$wsSummary.Range("E2").Value = "91%"
$wsSummary.Range("E2").Style = "Normal"

# Day 2 summary (17-05-2025)
$wsSummary.Range("A3").Value = "17-05-2025"
$wsSummary.Range("B3").Value = 10
$wsSummary.Range("C3").Value = 1
$wsSummary.Range("D3").Value = 11
$wsSummary.Range("E3").NumberFormat = "@"
$wsSummary.Range("E3").Value = "91%"
$wsSummary.Range("E3").Style = "Normal"

$wsResults.Select()
$excel.ActiveWindow.Zoom = 100
